$wb = $excel.ActiveWorkbook

# "numeric" sheet: add the "*" marker in B2 and move the selection to B3
$wsNumeric = $wb.Worksheets.Item("numeric")
$wsNumeric.Activate()
$wsNumeric.Range("B2").Value = "*"
[void]$wsNumeric.Range("B3").Select()

# "string" sheet: add the "*" marker in B2, move the selection to B3,
# and leave this as the active/selected sheet (it was "drop" before)
$wsString = $wb.Worksheets.Item("string")
$wsString.Activate()
$wsString.Range("B2").Value = "*"
[void]$wsString.Range("B3").Select()
